# Logged Week 17 data and fixed Simulate_Season.py tiebreaking method
$wb = $excel.ActiveWorkbook

# --- OFF sheet: update Home ("H") row totals ---
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B2").Value = 487
$wsOff.Range("C2").Value = 309
$wsOff.Range("D2").Value = 135
$wsOff.Range("E2").Value = 58

# --- DEF sheet: update Home ("H") row totals ---
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B2").Value = 493
$wsDef.Range("C2").Value = 332
$wsDef.Range("D2").Value = 135
$wsDef.Range("E2").Value = 59
$wsDef.Range("F2").Value = 14
